$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the existing data area (A1:C3) a distinct (but visually-equivalent)
# cell style so every populated cell ends up pointing at style index 1
# instead of the default style index 0.
$ws.Range("A1:C3").FormulaHidden = $true

# Row 4: present but blank (still carries the new style) - matches the
# target fixture's "blank but styled" row used to prove line-length isn't
# lost on XLSX -> CSV conversion.
$ws.Range("A4:C4").FormulaHidden = $true

# Row 5 is intentionally skipped/left absent.

# Row 6: single value in column A.
$ws.Range("A6").FormulaHidden = $true
$ws.Range("A6").Value = "another"

# Row 7: three values.
$ws.Range("A7:C7").FormulaHidden = $true
$ws.Range("A7").Value = "Extra"
$ws.Range("B7").Value = "line"
$ws.Range("C7").Value = "present"
